$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$qa = @(
    @("Nhà hàng Nét Huế Hà Nội chuyên phục vụ đặc sản của vùng miền nào?", "Nhà hàng Nét Huế Hà Nội chuyên phục vụ các món đặc sản Huế, mang đậm hương vị ẩm thực miền Trung như bún bò Huế, bánh khoái, bánh nậm,…"),
    @("Nhà Hàng Hibana By Koki ở Hà Nội có địa chỉ ở đâu?", "Nhà Hàng Hibana By Koki có địa chỉ tại số 11, đường Lê Phụng Hiểu, quận Hoàn Kiếm, Hà Nội."),
    @("Nhà hàng Buffet Chef Dzung ở Hà Nội nằm ở đâu?", "Nhà hàng Buffet Chef Dzung Hà Nội có địa chỉ tại 71 Đ. Nguyễn Chí Thanh, Thành Công, Đống Đa, Hà Nội."),
    @("Dimsum Corner Hà Nội nổi tiếng với điều gì và thực đơn của nhà hàng có những loại món nào?", "Dimsum Corner Hà Nội nổi tiếng là một trong những nhà hàng dimsum được nhiều tín đồ ẩm thực Trung Hoa yêu thích. Nhà hàng chinh phục thực khách nhờ thực đơn đa dạng và chất lượng món ăn ổn định. Điểm nổi bật nhất là thực đơn phong phú với hơn 30 món dimsum truyền thống như các loại há cảo, xíu mại, bánh bao và các món chiên hấp dẫn. Ngoài dimsum, nhà hàng còn phục vụ các món Trung Hoa khác như cháo, mì, cơm và các món tráng miệng, giúp bữa ăn trở nên trọn vẹn."),
    @("ManWah Taiwanese Hotpot có gì nổi bật khiến thực khách Hà Nội yêu thích và nhà hàng phục vụ những loại nước lẩu cũng như topping nào?", "ManWah Taiwanese Hotpot nổi bật bởi hương vị lẩu chuẩn Đài Loan, từ cay nồng đến thanh ngọt, đều rất đậm đà và cuốn hút. Tuy mới xuất hiện tại Việt Nam gần đây nhưng ManWah đã nhanh chóng thu hút thực khách nhờ hương vị đặc trưng này. Nhà hàng có đến 6 loại nước lẩu độc đáo phù hợp với nhiều sở thích khác nhau. Thực đơn của ManWah cũng rất đa dạng, với hơn 80 món topping, bao gồm các loại thịt bò, hải sản tươi ngon, cùng nhiều loại rau củ phong phú. Những yếu tố này mang đến một bữa lẩu đầy đủ và trọn vị, khiến ManWah trở thành lựa chọn yêu thích của nhiều tín đồ ẩm thực, đặc biệt trong những ngày gió mùa."),
    @("Lẩu Thái Deli Hà Nội có những cơ sở nào, giờ mở cửa ra sao và trải nghiệm ẩm thực mà quán mang lại được mô tả như thế nào?", "Lẩu Thái Deli Hà Nội có các cơ sở tại:`nNgõ 4 Hàm Nghi, Số 2 lô TT01, Moncity, Mỹ Đình`nB10, Biệt thự PANDORA, 53 Triều Khúc`n25 Vũ Ngọc Phan`nLK05 Biệt thự LACASA, Láng Hạ`nSố 5 Quan Hoa, Cầu Giấy`nGiờ mở cửa của Thái Deli là 11:00–14:00 và 17:00–23:00. Giá tham khảo ở đây nằm trong khoảng 100.000VNĐ – 160.000VNĐ/người.`nTrải nghiệm ẩm thực tại Thái Deli được mô tả là khiến thực khách vừa được chìm đắm trong không gian ấm cúng đậm nét văn hóa Thái, vừa thưởng thức hương vị ẩm thực chuẩn vị trứ danh Thái Lan. Khách đến đây sẽ có cảm giác thư giãn, thoải mái như đang lạc vào xứ sở Chùa Vàng, tách biệt khỏi náo nhiệt, đồng thời tận hưởng sự bùng nổ vị giác từ các món lẩu đặc trưng."),
    @("Lẩu cua đồng Song Hà ở Hà Nội có địa chỉ ở đâu, giờ mở cửa như thế nào, và món lẩu cua đồng tại đây được mô tả ra sao?", "Lẩu cua đồng Song Hà ở Hà Nội có địa chỉ tại 685 Lạc Long Quân, Quận Tây Hồ, Hà Nội. Giờ mở cửa của quán là 8:00 - 22:00, với mức giá tham khảo từ 50.000VNĐ - 300.000VNĐ.`nNước dùng lẩu cua đồng Song Hà được đảm bảo chế biến từ cua đồng 100%, vì vậy có vị ngọt thanh và mùi hương ngào ngạt đặc trưng của cua đồng. Khi nồi lẩu được mang ra, thực khách sẽ bị thu hút bởi nồi nước dùng nổi đầy gạch cua vàng óng, điểm thêm cà chua đỏ, đậu chiên vàng, cùng một chút hành lá và rau xanh. Để món lẩu cua đồng thêm trọn vẹn, quán phục vụ kèm bắp bò tươi, sườn sụn, nấm và các loại rau xanh khác, đầy đủ và hấp dẫn.")
)

$startRow = 62
for ($i = 0; $i -lt $qa.Length; $i++) {
    $r = $startRow + $i

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $qa[$i][0]
    $cellA.WrapText = $true

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $qa[$i][1]
    $cellB.WrapText = $true

    # Keep the pre-existing row height (autofit can change it when long
    # multi-line text is assigned).
    $ws.Rows.Item($r).RowHeight = 49.95
}

# Match the new selection / scroll position recorded in the saved file.
[void]$ws.Range("A69").Select()
